$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.916864633560181
$ws.Range("B1").Value = 2.781983613967896
$ws.Range("C1").Value = 3.343796968460083
$ws.Range("D1").Value = 1.097015261650085
$ws.Range("E1").Value = 0.7066972255706787
